# Scheduled-runner style refresh of the price/profit columns (H:N) on each
# of the eight Leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates per-row: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ).
# A few rows also gain/lose a sparsely-populated LeveProfit cell (M/N) to
# mirror upstream availability of NQ vs HQ market data for that item.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 393.625
$ws.Range("I5").Value = 435.57144
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 435.57144
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -320.57144
$ws.Range("N5").Value = -330

$ws.Range("H28").Value = 1168.5
$ws.Range("I28").Value = 343
$ws.Range("K28").Value = 343
$ws.Range("M28").Value = 142

$ws.Range("H57").Value = 72583.5
$ws.Range("J57").Value = 72583.5
$ws.Range("L57").Value = 217750.5
$ws.Range("N57").Value = -218748.5

$ws.Range("H80").Value = 6898.7144
$ws.Range("J80").Value = 7781.8335
$ws.Range("L80").Value = 23345.5005
$ws.Range("N80").Value = -25341.5005

$ws.Range("H83").Value = 6898.7144
$ws.Range("J83").Value = 7781.8335
$ws.Range("L83").Value = 70036.5015
$ws.Range("N83").Value = -80020.5015

$ws.Range("H86").Value = 8561.9375
$ws.Range("I86").Value = 10831.167
$ws.Range("J86").Value = 7200.4
$ws.Range("K86").Value = 10831.167
$ws.Range("L86").Value = 7200.4
$ws.Range("M86").Value = -9708.166999999999
$ws.Range("N86").Value = -9446.4

$ws.Range("H89").Value = 8561.9375
$ws.Range("I89").Value = 10831.167
$ws.Range("J89").Value = 7200.4
$ws.Range("K89").Value = 54155.835
$ws.Range("L89").Value = 36002
$ws.Range("M89").Value = -48539.835
$ws.Range("N89").Value = -47234

$ws.Range("H112").Value = 6496.524
$ws.Range("J112").Value = 7423.722
$ws.Range("L112").Value = 22271.166
$ws.Range("N112").Value = -24487.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1048212.25
$ws.Range("I2").Value = 1154976.8
$ws.Range("K2").Value = 1154976.8
$ws.Range("M2").Value = -1154863.8

$ws.Range("H32").Value = 18721.46
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 18721.46
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 18721.46
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -19295.46

$ws.Range("H110").Value = 1112405.9
$ws.Range("J110").Value = 1957.8334
$ws.Range("L110").Value = 1957.8334
$ws.Range("N110").Value = -6047.8334

$ws.Range("H116").Value = 1048212.25
$ws.Range("I116").Value = 1154976.8
$ws.Range("K116").Value = 1154976.8
$ws.Range("M116").Value = -1152682.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1048212.25
$ws.Range("I3").Value = 1154976.8
$ws.Range("K3").Value = 1154976.8
$ws.Range("M3").Value = -1154862.8

$ws.Range("H24").Value = 1000.3333
$ws.Range("I24").Value = 1000.3333
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1000.3333
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -765.3333
$ws.Range("N24").ClearContents()

$ws.Range("H86").Value = 3574831.2
$ws.Range("I86").Value = 5266175.5
$ws.Range("K86").Value = 5266175.5
$ws.Range("M86").Value = -5265052.5

$ws.Range("H89").Value = 3574831.2
$ws.Range("I89").Value = 5266175.5
$ws.Range("K89").Value = 26330877.5
$ws.Range("M89").Value = -26325261.5

$ws.Range("H94").Value = 5889391.5
$ws.Range("J94").Value = 11746.25
$ws.Range("L94").Value = 11746.25
$ws.Range("N94").Value = -12648.25

$ws.Range("H105").Value = 5210583
$ws.Range("I105").Value = 5684127
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 5684127
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = -5682380
$ws.Range("N105").Value = -5094

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 400.1111
$ws.Range("I22").Value = 387.6875
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 387.6875
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = -37.6875
$ws.Range("N22").Value = -1199.5

$ws.Range("H62").Value = 2916
$ws.Range("I62").Value = 3768.75
$ws.Range("J62").Value = 2347.5
$ws.Range("K62").Value = 3768.75
$ws.Range("L62").Value = 2347.5
$ws.Range("M62").Value = -3144.75
$ws.Range("N62").Value = -3595.5

$ws.Range("H65").Value = 2916
$ws.Range("I65").Value = 3768.75
$ws.Range("J65").Value = 2347.5
$ws.Range("K65").Value = 18843.75
$ws.Range("L65").Value = 11737.5
$ws.Range("M65").Value = -15723.75
$ws.Range("N65").Value = -17977.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5746.6562
$ws.Range("J81").Value = 6458.5356
$ws.Range("L81").Value = 19375.6068
$ws.Range("N81").Value = -21621.6068

$ws.Range("H84").Value = 5746.6562
$ws.Range("J84").Value = 6458.5356
$ws.Range("L84").Value = 58126.8204
$ws.Range("N84").Value = -69358.8204

$ws.Range("H119").Value = 24199.5
$ws.Range("I119").Value = 24199.5
$ws.Range("K119").Value = 72598.5
$ws.Range("M119").Value = -67760.5

$ws.Range("H137").Value = 1812.4117
$ws.Range("J137").Value = 1876.7
$ws.Range("L137").Value = 5630.1
$ws.Range("N137").Value = -15830.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 15900
$ws.Range("I24").Value = 7833.3335
$ws.Range("K24").Value = 7833.3335
$ws.Range("M24").Value = -7660.3335

$ws.Range("H80").Value = 24990198
$ws.Range("I80").Value = 43714850
$ws.Range("J80").Value = 23994.889
$ws.Range("K80").Value = 43714850
$ws.Range("L80").Value = 23994.889
$ws.Range("M80").Value = -43713852
$ws.Range("N80").Value = -25990.889

$ws.Range("H83").Value = 24990198
$ws.Range("I83").Value = 43714850
$ws.Range("J83").Value = 23994.889
$ws.Range("K83").Value = 218574250
$ws.Range("L83").Value = 119974.445
$ws.Range("M83").Value = -218569258
$ws.Range("N83").Value = -129958.445

$ws.Range("H113").Value = 5378856.5
$ws.Range("I113").Value = 6946273
$ws.Range("K113").Value = 6946273
$ws.Range("M113").Value = -6944103

$ws.Range("H132").Value = 1669.9103
$ws.Range("I132").Value = 1486.0145
$ws.Range("K132").Value = 4458.0435
$ws.Range("M132").Value = -1928.0435

$ws.Range("H139").Value = 99999.39999999999
$ws.Range("J139").Value = 99999.39999999999
$ws.Range("L139").Value = 99999.39999999999
$ws.Range("N139").Value = -110279.4

$ws.Range("H140").Value = 82598.60000000001
$ws.Range("J140").Value = 82598.60000000001
$ws.Range("L140").Value = 82598.60000000001
$ws.Range("N140").Value = -92958.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 47981250
$ws.Range("I82").Value = 65973840
$ws.Range("J82").Value = 1008.8333
$ws.Range("K82").Value = 65973840
$ws.Range("L82").Value = 1008.8333
$ws.Range("M82").Value = -65973479
$ws.Range("N82").Value = -1730.8333

$ws.Range("H85").Value = 47981250
$ws.Range("I85").Value = 65973840
$ws.Range("J85").Value = 1008.8333
$ws.Range("K85").Value = 65973840
$ws.Range("L85").Value = 1008.8333
$ws.Range("M85").Value = -65972592
$ws.Range("N85").Value = -3504.8333

$ws.Range("H132").Value = 6742.1304
$ws.Range("I132").Value = 6625.825
$ws.Range("K132").Value = 19877.475
$ws.Range("M132").Value = -17347.475

$ws.Range("H136").Value = 38228.406
$ws.Range("I136").Value = 52972.324
$ws.Range("K136").Value = 158916.972
$ws.Range("M136").Value = -156366.972

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6541.171
$ws.Range("I62").Value = 3029.8572
$ws.Range("J62").Value = 7264.0884
$ws.Range("K62").Value = 3029.8572
$ws.Range("L62").Value = 7264.0884
$ws.Range("M62").Value = -2405.8572
$ws.Range("N62").Value = -8512.088400000001

$ws.Range("H65").Value = 6541.171
$ws.Range("I65").Value = 3029.8572
$ws.Range("J65").Value = 7264.0884
$ws.Range("K65").Value = 15149.286
$ws.Range("L65").Value = 36320.442
$ws.Range("M65").Value = -12029.286
$ws.Range("N65").Value = -42560.442

$ws.Range("H81").Value = 7251755.5
$ws.Range("I81").Value = 11909728
$ws.Range("J81").Value = 6021
$ws.Range("K81").Value = 23819456
$ws.Range("L81").Value = 12042
$ws.Range("M81").Value = -23818395
$ws.Range("N81").Value = -14164

$ws.Range("H84").Value = 7251755.5
$ws.Range("I84").Value = 11909728
$ws.Range("J84").Value = 6021
$ws.Range("K84").Value = 119097280
$ws.Range("L84").Value = 60210
$ws.Range("M84").Value = -119091976
$ws.Range("N84").Value = -70818

$ws.Range("H96").Value = 3229.2307
$ws.Range("I96").Value = 2991.375
$ws.Range("J96").Value = 3609.8
$ws.Range("K96").Value = 2991.375
$ws.Range("L96").Value = 3609.8
$ws.Range("M96").Value = -1618.375
$ws.Range("N96").Value = -6355.8

$ws.Range("H126").Value = 2306.4119
$ws.Range("I126").Value = 2105.1667
$ws.Range("K126").Value = 6315.500100000001
$ws.Range("M126").Value = -3845.500100000001

$ws.Range("H136").Value = 836.3288
$ws.Range("I136").Value = 744.96826
$ws.Range("K136").Value = 2234.90478
$ws.Range("M136").Value = 315.0952200000002

$ws.Range("H138").Value = 86404.336
$ws.Range("J138").Value = 86404.336
$ws.Range("L138").Value = 86404.336
$ws.Range("N138").Value = -96684.336
